$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The worksheet is protected; temporarily unprotect to apply the data
# refresh, then restore protection afterwards.
$ws.Unprotect()

# Update the "as of" date in the confidential disclosure note (A41):
# 2021-05-26 -> 2021-05-27.
$ws.Range("A41").Value = "***CONFIDENTIAL***: For one-on-one client use only. Not approved for distribution.`nModel holdings provided as of 2021-05-27 for illustrative purposes only and are subject to change."

# Refresh the Weight (col D) and Percent Change (col E) figures for each
# holding row (2-38) with the latest model values.
$ws.Cells.Item(2,4).Value = 0.03105132460255736
$ws.Cells.Item(2,5).Value = -0.02074978204010458
$ws.Cells.Item(3,4).Value = 0.02941038836780022
$ws.Cells.Item(3,5).Value = -0.0171604402540535
$ws.Cells.Item(4,4).Value = 0.02878232321454884
$ws.Cells.Item(4,5).Value = 0.001021189685984192
$ws.Cells.Item(5,4).Value = 0.06313833792457726
$ws.Cells.Item(5,5).Value = -0.01073454287079334
$ws.Cells.Item(6,4).Value = 0.01520234797361972
$ws.Cells.Item(6,5).Value = 0.0139662672670382
$ws.Cells.Item(7,4).Value = 0.01544135306272646
$ws.Cells.Item(7,5).Value = 0.004132541888947294
$ws.Cells.Item(8,4).Value = 0.02767702136233675
$ws.Cells.Item(8,5).Value = -0.01358205826870684
$ws.Cells.Item(9,4).Value = 0.03562645438270749
$ws.Cells.Item(9,5).Value = 0.0003256621797653736
$ws.Cells.Item(10,4).Value = 0.0291423778067145
$ws.Cells.Item(10,5).Value = 0.002985906521219883
$ws.Cells.Item(11,4).Value = 0.02866630132663294
$ws.Cells.Item(11,5).Value = 0.00453300594957029
$ws.Cells.Item(12,4).Value = 0.01124638833531408
$ws.Cells.Item(12,5).Value = 0.003008940852819997
$ws.Cells.Item(13,4).Value = 0.01433721142939353
$ws.Cells.Item(13,5).Value = -0.0003776435045318438
$ws.Cells.Item(14,4).Value = 0.01415853772200305
$ws.Cells.Item(14,5).Value = -0.006214149139579295
$ws.Cells.Item(15,4).Value = 0.00893059145251282
$ws.Cells.Item(15,5).Value = 0.02217218084185002
$ws.Cells.Item(16,4).Value = 0.008136034890102135
$ws.Cells.Item(16,5).Value = 0.02032085561497321
$ws.Cells.Item(17,4).Value = 0.03108922508594322
$ws.Cells.Item(17,5).Value = -0.01676245210727956
$ws.Cells.Item(18,4).Value = 0.0244251212138674
$ws.Cells.Item(18,5).Value = -0.01001480449359915
$ws.Cells.Item(19,4).Value = 0.03322866869911232
$ws.Cells.Item(19,5).Value = 0.001134776536312776
$ws.Cells.Item(20,4).Value = 0.03167977649543512
$ws.Cells.Item(20,5).Value = 0.01553439540987589
$ws.Cells.Item(21,4).Value = 0.04705712414999465
$ws.Cells.Item(21,5).Value = -0.01274691497536495
$ws.Cells.Item(22,4).Value = 0.03552203468358318
$ws.Cells.Item(22,5).Value = 0.01080566140446382
$ws.Cells.Item(23,4).Value = 0.03068102141029246
$ws.Cells.Item(23,5).Value = 0.004537862792676428
$ws.Cells.Item(24,4).Value = 0.02948812303270387
$ws.Cells.Item(24,5).Value = 0.01664305949008482
$ws.Cells.Item(25,4).Value = 0.01609629662001169
$ws.Cells.Item(25,5).Value = 0.02454319385879566
$ws.Cells.Item(26,4).Value = 0.01536709905446029
$ws.Cells.Item(26,5).Value = 0.00166100415251047
$ws.Cells.Item(27,4).Value = 0.03049790019719853
$ws.Cells.Item(27,5).Value = 0.003068768307992853
$ws.Cells.Item(28,4).Value = 0.03005740376274452
$ws.Cells.Item(28,5).Value = 0.00005146680391132641
$ws.Cells.Item(29,4).Value = 0.02917834459196843
$ws.Cells.Item(29,5).Value = -0.00866833671318934
$ws.Cells.Item(30,4).Value = 0.02911743310081258
$ws.Cells.Item(30,5).Value = -0.002410694718386952
$ws.Cells.Item(31,4).Value = 0.03643087280559102
$ws.Cells.Item(31,5).Value = -0.01350318471337586
$ws.Cells.Item(32,4).Value = 0.03154441762619992
$ws.Cells.Item(32,5).Value = -0.01207012811867858
$ws.Cells.Item(33,4).Value = 0.02916016782952827
$ws.Cells.Item(33,5).Value = -0.00198275862068964
$ws.Cells.Item(34,4).Value = 0.03196093620381797
$ws.Cells.Item(34,5).Value = -0.001234239248808056
$ws.Cells.Item(35,4).Value = 0.03004502809470016
$ws.Cells.Item(35,5).Value = 0
$ws.Cells.Item(36,4).Value = 0.03129303686904881
$ws.Cells.Item(36,5).Value = 0.01557189643452994
$ws.Cells.Item(37,4).Value = 0.03513297461943855
$ws.Cells.Item(37,5).Value = 0.001144819690898746
$ws.Cells.Item(38,5).Value = -0.001436544342211854

$ws.Protect()
